# Bugfix in Edge Weighting Script
# - Append a new edge row (TP 500 Team Project Split 2 -> Scientific Work (SW))
# - Apply an AutoFilter on column "Outgoing" (A) limited to "TP 500 Team Project",
#   which hides every other data row except the header and the new appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the new data row right after the last existing data row (row 164 -> 165)
$newRow = 165
$ws.Cells.Item($newRow, 1).Value = "TP 500 Team Project Split 2"
$ws.Cells.Item($newRow, 2).Value = "Scientific Work (SW)"

# Apply the AutoFilter on the original data range (A1:B164), filtering column A
# ("Outgoing") down to the single value "TP 500 Team Project" (row 57).
$filterRange = $ws.Range("A1:B164")
$filterRange.AutoFilter(1, @("TP 500 Team Project"), 7)

# Move selection to the new (now-visible, below-filter) row, mirroring the
# author's final cursor position after the edit.
$ws.Range("A166").Select()
